$wb = $excel.ActiveWorkbook

# Rename sheets
$wb.Worksheets.Item(1).Name = "Export as TSV"
$wb.Worksheets.Item(5).Name = "gdna_fragmenta...assurance list"

$ws = $wb.Worksheets.Item(1)

# Freeze top row (pane split) - select A2 first so the engine infers a pure
# frozen state (not frozenSplit) with topLeftCell A2 / activePane bottomLeft
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true

# Update the data validation that references the renamed list sheet
$ws.Range("O2:O1048576").Validation.Modify(3, 1, 1, "'gdna_fragmenta...assurance list'!`$A`$1:`$A`$2")

# Add errorTitle / error messages to each data validation
$ws.Range("I2:I1048576").Validation.ErrorTitle = "Value must come from list"
$ws.Range("I2:I1048576").Validation.ErrorMessage = "Value must be one of: sequence."

$ws.Range("J2:J1048576").Validation.ErrorTitle = "Value must come from list"
$ws.Range("J2:J1048576").Validation.ErrorMessage = "Value must be one of: WGS."

$ws.Range("K2:K1048576").Validation.ErrorTitle = "Value must come from list"
$ws.Range("K2:K1048576").Validation.ErrorMessage = "Value must be one of: DNA."

$ws.Range("L2:L1048576").Validation.ErrorTitle = "Not a boolean"
$ws.Range("L2:L1048576").Validation.ErrorMessage = 'The values in this column must be "TRUE" or "FALSE".'

$ws.Range("O2:O1048576").Validation.ErrorTitle = "Value must come from list"
$ws.Range("O2:O1048576").Validation.ErrorMessage = "Value must be one of: Pass / Fail."

$ws.Range("P2:P1048576").Validation.ErrorTitle = "Not a number"
$ws.Range("P2:P1048576").Validation.ErrorMessage = "The values in this column must be numbers."

$ws.Range("Q2:Q1048576").Validation.ErrorTitle = "Value must come from list"
$ws.Range("Q2:Q1048576").Validation.ErrorMessage = "Value must be one of: ug."

$ws.Range("T2:T1048576").Validation.ErrorTitle = "Value must come from list"
$ws.Range("T2:T1048576").Validation.ErrorMessage = "Value must be one of: single-end / paired-end."

$ws.Range("V2:V1048576").Validation.ErrorTitle = "Not a number"
$ws.Range("V2:V1048576").Validation.ErrorMessage = "The values in this column must be numbers."

$ws.Range("W2:W1048576").Validation.ErrorTitle = "Value must come from list"
$ws.Range("W2:W1048576").Validation.ErrorMessage = "Value must be one of: ng."

$ws.Range("X2:X1048576").Validation.ErrorTitle = "Not a number"
$ws.Range("X2:X1048576").Validation.ErrorMessage = "The values in this column must be numbers."

$ws.Range("AA2:AA1048576").Validation.ErrorTitle = "Not a number"
$ws.Range("AA2:AA1048576").Validation.ErrorMessage = "The values in this column must be numbers."

$ws.Range("AB2:AB1048576").Validation.ErrorTitle = "Not a number"
$ws.Range("AB2:AB1048576").Validation.ErrorMessage = "The values in this column must be numbers."

Write-Host "Done"
